$wb = $excel.ActiveWorkbook

function Set-Cells {
    param(
        [string]$SheetName,
        [hashtable]$Updates
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($ref in $Updates.Keys) {
        $ws.Range($ref).Value = $Updates[$ref]
    }
}

# Sheet "展览" (sheetId 1)
Set-Cells "展览" @{
    "F3"  = 3479
    "F5"  = 8148
    "F7"  = 71
    "F8"  = 2128
    "G9"  = 68
    "F12" = 1092
    "F16" = 1163
    "F18" = 9
    "F19" = 734
    "F20" = 508
    "F21" = 522
    "F22" = 61
    "F24" = 6784
    "F25" = 114
    "F26" = 53620
    "F27" = 4148
    "F29" = 783
    "F30" = 362
    "F31" = 71
    "F34" = 579
    "F35" = 1621
    "F36" = 566
    "F38" = 833
    "F39" = 1050
    "F40" = 379
    "F44" = 683
}

# Sheet "演出" (sheetId 2)
Set-Cells "演出" @{
    "F12" = 40
    "F13" = 94
    "F17" = 7351
    "F28" = 109
    "F30" = 76
    "F37" = 20
}

# Sheet "本地生活" (sheetId 3)
Set-Cells "本地生活" @{
    "F4"  = 2240
    "F5"  = 1499
    "F7"  = 637
    "F9"  = 9275
    "F10" = 1553
    "F15" = 94
}

# Sheet "全部类型" (sheetId 4)
Set-Cells "全部类型" @{
    "F3"  = 3479
    "F4"  = 2240
    "F6"  = 637
    "F7"  = 1553
    "F10" = 71
    "F14" = 1163
    "F15" = 94
    "F16" = 9
    "F17" = 734
    "F18" = 61
    "F20" = 6784
    "F21" = 114
    "F22" = 53620
    "F27" = 4148
    "F29" = 362
    "F30" = 71
    "F32" = 579
    "F33" = 94
    "F35" = 566
    "F38" = 833
    "F50" = 20
}
